$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '65.338.74'
$ws.Range('E2').Value = '  +3.42%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.492.30'
$ws.Range('E3').Value = '  +2.89%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '581.66'
$ws.Range('E5').Value = '  +2.61%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '163.54'
$ws.Range('E6').Value = '  +5.44%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.611'
$ws.Range('E7').Value = '  +12.46%  '
$ws.Range('E8').Value = '  -0.02%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '3.494.52'
$ws.Range('E9').Value = '  +3.04%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '7.28'
$ws.Range('E10').Value = '  -1.80%  '
$ws.Range('E11').Value = '  +3.63%  '
$ws.Range('E12').Value = '  +3.89%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.093.20'
$ws.Range('E13').Value = '  +2.84%  '
$ws.Range('E14').Value = '  +0.46%  '
$ws.Range('E15').Value = '  +1.81%  '
$ws.Range('E16').Value = '  +6.17%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '65.313.70'
$ws.Range('E17').Value = '  +3.31%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '3.476.44'
$ws.Range('E18').Value = '  +2.06%  '
$ws.Range('E19').Value = '  +3.83%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '14.43'
$ws.Range('E20').Value = '  +2.41%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '384.85'
$ws.Range('E21').Value = '  +1.75%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '8.26'
$ws.Range('E22').Value = '  +2.63%  '
$ws.Range('E23').Value = '  +5.06%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '72.70'
$ws.Range('E25').Value = '  +0.28%  '
$ws.Range('E26').Value = '  +1.59%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '10.09'
$ws.Range('E27').Value = '  +7.27%  '
$ws.Range('E28').Value = '  +0.56%  '
$ws.Range('B29').Value = 'Fetch.AI'
$ws.Range('C29').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.55'
$ws.Range('E29').Value = '  +14.02%  '
$ws.Range('B30').Value = 'Binance-PegBSC-USD'
$ws.Range('C30').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.00'
$ws.Range('E30').Value = '  -0.02%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '6.21'
$ws.Range('E31').Value = '  +2.41%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '2.07'
$ws.Range('E32').Value = '  +3.35%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '23.78'
$ws.Range('E33').Value = '  +2.75%  '
$ws.Range('E34').Value = '  +6.35%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.64'
$ws.Range('E35').Value = '  +13.31%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '162.55'
$ws.Range('E36').Value = '  +1.70%  '
$ws.Range('E37').Value = '  +6.15%  '
$ws.Range('B38').Value = 'Hedera'
$ws.Range('C38').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.0784'
$ws.Range('E38').Value = '  +4.23%  '
$ws.Range('B39').Value = 'Maker'
$ws.Range('C39').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '3.014.41'
$ws.Range('E39').Value = '  +2.00%  '
$ws.Range('B40').Value = 'RenderToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '6.85'
$ws.Range('E40').Value = '  +8.24%  '
$ws.Range('B41').Value = 'EnergySwap'
$ws.Range('C41').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '27.02'
$ws.Range('E41').Value = '  -0.01%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '4.60'
$ws.Range('E42').Value = '  +6.50%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.0322'
$ws.Range('E43').Value = '  +1.89%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '43.14'
$ws.Range('E44').Value = '  +3.35%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.783'
$ws.Range('E45').Value = '  +3.01%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '26.03'
$ws.Range('E46').Value = '  +12.07%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.13'
$ws.Range('E47').Value = '  +5.01%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '323.51'
$ws.Range('E48').Value = '  +11.22%  '
$ws.Range('B49').Value = 'SuiNetwork'
$ws.Range('C49').Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.889'
$ws.Range('E49').Value = '  +7.08%  '
$ws.Range('B50').Value = 'Cosmos'
$ws.Range('C50').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '6.76'
$ws.Range('E50').Value = '  +6.55%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.109'
$ws.Range('E51').Value = '  +6.54%  '
